$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp caption (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 17:52"

# --- Update per-country statistics (columns: B=Casos totales, C=Nuevos casos,
#     D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 800639
$ws.Range("C4").Value = 7880
$ws.Range("D4").Value = 72898
$ws.Range("E4").Value = 684548
$ws.Range("G4").Value = 679
$ws.Range("H4").Value = 43193

# Alemania (row 8)
$ws.Range("B8").Value = 147786
$ws.Range("C8").Value = 721
$ws.Range("E8").Value = 47687
$ws.Range("G8").Value = 37
$ws.Range("H8").Value = 4899

# Canada (row 16)
$ws.Range("B16").Value = 37398
$ws.Range("C16").Value = 569
$ws.Range("E16").Value = 23084
$ws.Range("G16").Value = 38
$ws.Range("H16").Value = 1728

# Chile (row 28)
$ws.Range("F28").Value = 392

# Polonia (row 31)
$ws.Range("B31").Value = 9856
$ws.Range("C31").Value = 263
$ws.Range("E31").Value = 8158
$ws.Range("G31").Value = 21
$ws.Range("H31").Value = 401

# Rumania (row 33)
$ws.Range("E33").Value = 6591
$ws.Range("G33").Value = 20
$ws.Range("H33").Value = 498

# Dinamarca (row 37)
$ws.Range("F37").Value = 81

# Chequia (row 40)
$ws.Range("B40").Value = 6961
$ws.Range("C40").Value = 61
$ws.Range("D40").Value = 1753
$ws.Range("E40").Value = 5007
$ws.Range("G40").Value = 7
$ws.Range("H40").Value = 201

# Argentina (row 57)
$ws.Range("E57").Value = 2044
$ws.Range("G57").Value = 5
$ws.Range("H57").Value = 147

# Grecia (row 61)
$ws.Range("B61").Value = 2401
$ws.Range("C61").Value = 156
$ws.Range("E61").Value = 1703
$ws.Range("F61").Value = 59
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 121

# Rows 81-82: Cuba overtakes Afganistan in ranking
$ws.Range("A81").Value = "Cuba"
$ws.Range("B81").Value = 1137
$ws.Range("C81").Value = 50
$ws.Range("D81").Value = 309
$ws.Range("E81").Value = 790
$ws.Range("F81").Value = 18
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 38

$ws.Range("A82").Value = "Afganistan"
$ws.Range("B82").Value = 1092
$ws.Range("C82").Value = 66
$ws.Range("D82").Value = 150
$ws.Range("E82").Value = 906
$ws.Range("F82").Value = 7
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 36

# Rows 86-88: Republica de Yibuti overtakes Tunez and Costa de Marfil in ranking
$ws.Range("A86").Value = "Republica de Yibuti"
$ws.Range("B86").Value = 945
$ws.Range("C86").Value = 99
$ws.Range("D86").Value = 112
$ws.Range("E86").Value = 831
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 2

$ws.Range("A87").Value = "Tunez"
$ws.Range("B87").Value = 884
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 148
$ws.Range("E87").Value = 698
$ws.Range("F87").Value = 34
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 38

$ws.Range("A88").Value = "Costa de Marfil"
$ws.Range("B88").Value = 879
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 287
$ws.Range("E88").Value = 582
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 10

# Mauricio (row 113)
$ws.Range("D113").Value = 243
$ws.Range("E113").Value = 76
